$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing user from Jason Poage -> Rachel Tipton, with fuller data ---
$ws.Range("A2").Value = "Rachel"
$ws.Range("B2").Value = "Tipton"
$ws.Range("C2").Value = "Rachel.Tipton"
$ws.Range("D2").Value = "Rachel.Tipton@simplymac.com"

# E2 ("316") must stay text, like the rest of the sheet's inline strings -
# force text format before assigning so Excel doesn't auto-coerce it to a
# number, then drop back to the default "Normal" style so no stray
# formatting is left on the cell.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "316"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").Value = "Operation Manager"
$ws.Range("G2").Value = "316 Louisville:23.126.11.33"
$ws.Range("H2").Value = "UT101 - Orem:regional_manager|UT102 - Fort Union:regional_manager|UT105 - St. George:regional_manager|UT106 - Foothill:regional_manager|ID401 - Idaho Falls:regional_manager|MO510 - Springfield MO:regional_manager|TX203 - Tyler:regional_manager|OR410 - Bend:regional_manager|TX200 - Lubbock:regional_manager|MT404 - Missoula:regional_manager|MT406 - Bozeman:regional_manager|UT100 - Park City:regional_manager|AR205 - Rogers:regional_manager|TN312 - Nashville:regional_manager|TN313 - Vanderbilt:regional_manager|TN314 - Chattanooga:regional_manager|TN315 - Johnson City:regional_manager|KY316 - Louisville:inventory_manager,shift_lead,technical_manager,sales_manager,regional_manager|IN305 - Ft. Wayne:regional_manager|IN301 - Evansville:regional_manager|GA702 - Athens:regional_manager|VA704 - Charlottesville:regional_manager|GA705 - Macon:regional_manager|GA706 - West Cobb:regional_manager|GA707 - Peachtree City:regional_manager|GA708 - Norcross:regional_manager|GA710 - Forsyth:regional_manager|GA711 - East Cobb:regional_manager|AL712 - Montgomery:regional_manager|GA713 - St. Simons:regional_manager|OR419 - Corvallis:regional_manager|OR420 - Eugene:regional_manager|OH308 - Cincinnati:regional_manager|TX206 - Katy:regional_manager|FL917 - Miami Warehouse:regional_manager|FL714 - Orlando:regional_manager|FL715 - Sawgrass:regional_manager|FL716 - Dolphin:regional_manager|NC317 - Asheville:regional_manager|GA717 - Gainesville:regional_manager|OR421 - Hillsboro:regional_manager|FL801 - Destin:regional_manager|TX207 - Waco:regional_manager|MO514 - Columbia:regional_manager|FL805 - Orlando UCF:regional_manager|TX208 - San Marcos:regional_manager|FL802 - Tallahassee:regional_manager|NC605 - Winston Salem:regional_manager|SC604 - Myrtle Beach:regional_manager|AL803 - Tuscaloosa:regional_manager|KS511 - Lawrence:regional_manager|SC602 - Columbia SC:regional_manager|NC603 - Wilmington:regional_manager|SC601 - Mt. Pleasant:regional_manager|KS513 - Wichita:regional_manager|FL807 - Clearwater:regional_manager|FL804 - Jacksonville:regional_manager"
$ws.Range("I2").Value = "yes"
$ws.Range("J2").Value = "User"
$ws.Range("K2").Value = "Email"
$ws.Range("L2").Value = "no"
$ws.Range("M2").Value = "None"
$ws.Range("N2").Value = "None"

# --- Row 3: new user Dirk Tomlinson ---
$ws.Range("A3").Value = "Dirk"
$ws.Range("B3").Value = "Tomlinson"
$ws.Range("C3").Value = "Dirk.Tomlinson"
$ws.Range("D3").Value = "Dirk.Tomlinson@simplymac.com"

# E3 ("301") - same text-coercion concern as E2 above.
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "301"
$ws.Range("E3").Style = "Normal"

$ws.Range("F3").Value = "Technician"
$ws.Range("G3").Value = "301 Evansville:23.28.217.91"
$ws.Range("H3").Value = "IN301 - Evansville:technical_manager,sales_manager"
$ws.Range("I3").Value = "yes"
$ws.Range("J3").Value = "User"
$ws.Range("K3").Value = "Email"
$ws.Range("L3").Value = "no"
$ws.Range("M3").Value = "None"
$ws.Range("N3").Value = "None"
